$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    7  = -6
    10 = 1
    14 = -1
    16 = -4
    17 = -1
    19 = -1
    21 = -7
    31 = -7
    33 = -4
    37 = -9
    38 = 3
    40 = 4
    43 = -6
    46 = -2
    47 = -7
    48 = 1
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
